$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value from "Emre Abale" to "Rob Oudman"
$ws.Range("A2").Value = "Rob Oudman"

# Update selection to A2 (was B4)
$ws.Range("A2").Select()
